# Apply updated attribution values to row 2 (relative direction update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "A2" = -0
    "B2" = -0.0736143064681548
    "C2" = -0
    "D2" = 0.2023283625086515
    "E2" = 0.005361901149070607
    "G2" = 0
    "I2" = -0
    "J2" = -0
    "K2" = 0.005910374655943606
    "L2" = -0
    "M2" = 0.2017405004068997
    "N2" = -0.003175673222564392
    "R2" = -0
    "S2" = 0
    "T2" = -0.0836551046379089
    "V2" = 0.01416007321150566
    "W2" = -0.03298341659304817
    "Y2" = -0
    "Z2" = -0
    "AB2" = 0
    "AC2" = -0.05416417955287071
    "AD2" = 0
    "AE2" = -0.01788162495550331
    "AF2" = 0.0002069792777307436
    "AG2" = -0
    "AH2" = -0
    "AI2" = -0
    "AJ2" = 0
    "AK2" = -0
    "AL2" = -0.03184932082569965
    "AM2" = 0
    "AN2" = 0.02827770634814052
    "AO2" = 0.06944358562979185
    "AQ2" = 0
    "AR2" = -0
    "AT2" = 0
    "AU2" = -0.1497027310705481
    "AW2" = 0.07080831603100772
    "AX2" = -0.001817919973327277
    "AY2" = -0
    "BC2" = -0
    "BD2" = -0.01418805710578807
    "BF2" = 0.08649163433815991
    "BG2" = 0.03283123518905573
    "BJ2" = -0
    "BL2" = 0
    "BM2" = 0.03190481457958391
    "BO2" = -0.04257541708426302
    "BP2" = -0.08905310676590357
    "BU2" = 0
    "BV2" = -0.04640710802875297
    "BW2" = 0
    "BX2" = 0.01094127294829059
    "BY2" = -0.02004983166574863
    "BZ2" = -0
    "CB2" = 0
    "CD2" = -0
    "CE2" = 0.03274941356648393
    "CG2" = -0.03099117283538349
    "CH2" = 0.01593166186343111
    "CJ2" = -0
    "CM2" = -0
    "CN2" = -0.01061381960660221
    "CP2" = 0.02133540246658532
    "CQ2" = 0.03716094318380431
    "CT2" = 0
    "CU2" = -0
    "CV2" = -0
    "CW2" = 0.04597314766486385
    "CY2" = -0.03342648399499332
    "CZ2" = 0.01017704690408558
    "DD2" = -0
    "DE2" = -0
    "DF2" = 0.02874569132567836
    "DH2" = 0.02900526664094873
    "DI2" = 0.03380050877759293
    "DJ2" = 0
    "DK2" = -0
    "DL2" = -0
    "DN2" = 0
    "DO2" = -0.01950247745448723
    "DP2" = -0
    "DQ2" = 0.03637034262361485
    "DR2" = -0.01945341551444906
    "DS2" = -0
    "DW2" = 0
    "DX2" = -0.05668836815106189
    "DY2" = -0
    "DZ2" = -0.008248668484950638
    "EA2" = -0.02434073422596091
    "EB2" = 0
    "EF2" = -0
    "EG2" = 0.04085971144248264
    "EI2" = 0.06686095049629477
    "EJ2" = -0.02407360759003618
    "EO2" = 0
    "EP2" = 0.04575396385905522
    "EQ2" = 0
    "ER2" = -0.0349068518066118
    "ES2" = 0.03654119765287879
    "ET2" = 0
    "EU2" = -0
    "EV2" = 0
    "EX2" = 0
    "EY2" = 0.04340351386436194
    "FA2" = -0.02669773779825179
    "FB2" = 0.01710822948871973
    "FD2" = -0
    "FG2" = -0
    "FH2" = 0.001691558334483329
    "FI2" = 0
    "FJ2" = -0.006288705109680439
    "FK2" = -0.008476023910300627
    "FL2" = -0
    "FN2" = -0
    "FP2" = -0
    "FQ2" = -0.01406165917757605
    "FR2" = -0
    "FS2" = -0.0184191401413391
    "FT2" = 0.005010644615690384
    "FV2" = -0
    "FW2" = -0
    "FY2" = 0
    "FZ2" = -0.03040407898268115
    "GB2" = 0.03399166782388464
    "GD2" = 0
    "GE2" = -0
}

foreach ($addr in $newValues.Keys) {
    $ws.Range($addr).Value = $newValues[$addr]
}

Write-Output "Applied $($newValues.Count) cell updates to row 2"
